$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the remaining cells for row 7
$ws.Range("X7").Value = 0.47000100000001055
$ws.Range("Y7").Value = "Up"

# Add new row 8
$ws.Range("A8").Value = 42649.886747685188
$ws.Range("B8").Value = 6
$ws.Range("C8").Value = "Buy"
$ws.Range("D8").Value = 36
$ws.Range("E8").Value = 12329
$ws.Range("F8").Value = 2332
$ws.Range("G8").Value = 54
$ws.Range("H8").Value = 42
$ws.Range("I8").Value = 85
$ws.Range("J8").Value = 12
$ws.Range("K8").Value = 15419
$ws.Range("L8").Value = 238
$ws.Range("M8").Value = 183
$ws.Range("N8").Value = 113
$ws.Range("O8").Value = 17
$ws.Range("P8").Value = "Noun"
$ws.Range("Q8").Value = 35.958706302092025
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = -0.0019
$ws.Range("T8").Value = -0.023
$ws.Range("U8").Value = 14.71
$ws.Range("V8").Value = "N/A"
$ws.Range("W8").Value = 0

# Apply same number formats as the row above (date format on A, percent on S/T)
# by copying formats only, so the existing style indexes are reused instead
# of creating brand-new custom number formats.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)

$ws.Range("S7").Copy()
$ws.Range("S8").PasteSpecial(-4122)

$ws.Range("T7").Copy()
$ws.Range("T8").PasteSpecial(-4122)

$excel.CutCopyMode = 0
